$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows after the header row (new rows 2 through 8),
# shifting the existing data rows down (old row 2 becomes row 9, etc.)
$ws.Rows("2:8").Insert()
$ws.Range("A2:C8").ClearFormats()

$topData = New-Object 'object[,]' 7,3
$topData[0,0] = 1.19674801826477
$topData[0,1] = 1.663910716772079
$topData[0,2] = 2.702408194541931
$topData[1,0] = -0.4646213054656981
$topData[1,1] = 2.042550325393677
$topData[1,2] = 2.406269252300262
$topData[2,0] = -0.3934619426727292
$topData[2,1] = 1.991465017199517
$topData[2,2] = 1.958218067884445
$topData[3,0] = -0.0652618408203127
$topData[3,1] = 1.843156695365906
$topData[3,2] = 2.04642915725708
$topData[4,0] = -0.2364732027053833
$topData[4,1] = 1.819270551204681
$topData[4,2] = 2.093152940273285
$topData[5,0] = -0.3396859169006348
$topData[5,1] = 1.832332909107209
$topData[5,2] = 2.310090780258179
$topData[6,0] = -0.2427999973297116
$topData[6,1] = 1.836586102843285
$topData[6,2] = 2.258781224489212
$ws.Range("A2:C8").Value = $topData

# Append 3 new rows at the bottom of the data (rows 29-31)
$bottomData = New-Object 'object[,]' 3,3
$bottomData[0,0] = -0.2655735015869125
$bottomData[0,1] = 2.233672142028808
$bottomData[0,2] = 0.9439086914062514
$bottomData[1,0] = 0.07992589473724532
$bottomData[1,1] = 1.95888604223728
$bottomData[1,2] = 1.25704461336136
$bottomData[2,0] = 0.3567421436309829
$bottomData[2,1] = 2.357963830232623
$bottomData[2,2] = 1.160924613475799
$ws.Range("A29:C31").Value = $bottomData
